# Generate Report for Handoff
# Updates the localization-status workbook:
#  - renames the tracked markdown file (new GUID-based name)
#  - adds a brand new source file row ("ffff62d6...md")
#  - flips status from "Handoff transform failed" -> "Ready for handoff"
#  - records the new handoff (.xlf) artifacts + timestamps for zh-cn / de-de
#  - flips the handoff reason from "Ignored" -> "Include" for the two
#    real source files, leaving ".localization-config" as "Ignored"

$wb = $excel.ActiveWorkbook

# ---- shared literal values -------------------------------------------------
$oldMd       = "5758bb67-cef7-48ca-8c47-cf2be1947d43.md"
$md1         = "48d6ee85-c3e0-4bc4-9566-5b30de1e9115.md"
$md2         = "ffff62d6325c-0ff1-49c1-8d54-4407761469fa.md"
$cfgName     = ".localization-config"

$xlfZh       = "48d6ee85-c3e0-4bc4-9566-5b30de1e9115.bb477888bd80f5db979b04e7359a08a48da61915.zh-cn.xlf"
$xlfDe       = "48d6ee85-c3e0-4bc4-9566-5b30de1e9115.bb477888bd80f5db979b04e7359a08a48da61915.de-de.xlf"
$tsZh        = "2016-02-18 10:17:43"
$tsDe        = "2016-02-18 10:17:54"

$statusReady = "Ready for handoff"
$statusNotLoc= "Not to be localized"
$reasonInc   = "Include"
$reasonIgn   = "Ignored"
$epoch       = "0001-01-01 00:00:00"
$dateFmt     = "yyyy-mm-dd HH:mm:ss"

$hlColor     = 15570276  # OLE BGR packing of RGB FF6495ED (matches workbook's HyperLink style)

$mdUrlBase   = "https://github.com/OpenLocalizationTest/oltest/blob/66af2ca68e3bbbd78f84de0bdb488713059e4662/e2e/"
$cfgUrl      = "https://github.com/OpenLocalizationTest/oltest/blob/b8dc8fe64cf5c43d454d364756b10a25ffb8e9e5/.localization-config"
$xlfUrlBase  = "https://github.com/OpenLocalizationTest/oltest/blob/bb477888bd80f5db979b04e7359a08a48da61915/"

function Style-Hyperlink($range) {
    $range.Font.Underline = 2
    $range.Font.Color = $hlColor
}

# =============================================================================
# Sheet "Overview"
# =============================================================================
$ws = $wb.Worksheets.Item("Overview")

# Make room for the new source-file row (old row 3 -> row 4).
$ws.Rows.Item(3).Insert()

# Clear all existing hyperlinks on the sheet; we'll recreate them below so
# none of them keep stale anchors after the row shift.
$ws.Hyperlinks.Delete()

# Row 2 - renamed markdown file, status flips to "Ready for handoff"
$ws.Range("B2").Value2 = $statusReady
$ws.Range("C2").Value2 = $statusReady

# Row 3 - brand new source file
$ws.Range("B3").Value2 = $statusReady
$ws.Range("C3").Value2 = $statusReady

# Row 4 - old ".localization-config" row, shifted down, content unchanged
$ws.Range("B4").Value2 = $statusNotLoc
$ws.Range("C4").Value2 = $statusNotLoc

$ws.Hyperlinks.Add($ws.Range("A2"), ($mdUrlBase + $md1), "", "", $md1) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), ($mdUrlBase + $md2), "", "", $md2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), $cfgUrl, "", "", $cfgName) | Out-Null

Style-Hyperlink $ws.Range("A2")
Style-Hyperlink $ws.Range("A3")
Style-Hyperlink $ws.Range("A4")

# =============================================================================
# Sheet "zh-cn"
# =============================================================================
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Rows.Item(3).Insert()
$ws.Hyperlinks.Delete()

# Row 2
$ws.Range("B2").Value2 = $statusReady
$ws.Range("D2").Value2 = $tsZh
$ws.Range("D2").NumberFormat = $dateFmt
$ws.Range("G2").Value2 = $epoch
$ws.Range("H2").Value2 = $reasonInc

# Row 3 (new)
$ws.Range("B3").Value2 = $statusReady
$ws.Range("D3").Value2 = $tsZh
$ws.Range("D3").NumberFormat = $dateFmt
$ws.Range("G3").Value2 = $epoch
$ws.Range("H3").Value2 = $reasonInc

# Row 4 (shifted ".localization-config" row)
$ws.Range("B4").Value2 = $statusNotLoc
$ws.Range("D4").Value2 = $epoch
$ws.Range("D4").NumberFormat = $dateFmt
$ws.Range("G4").Value2 = $epoch
$ws.Range("H4").Value2 = $reasonIgn

$ws.Hyperlinks.Add($ws.Range("A2"), ($mdUrlBase + $md1), "", "", $md1) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), ($xlfUrlBase + $xlfZh), "", "", $xlfZh) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), ($mdUrlBase + $md2), "", "", $md2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), ($xlfUrlBase + $xlfZh), "", "", $xlfZh) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), $cfgUrl, "", "", $cfgName) | Out-Null

Style-Hyperlink $ws.Range("A2")
Style-Hyperlink $ws.Range("C2")
Style-Hyperlink $ws.Range("A3")
Style-Hyperlink $ws.Range("C3")
Style-Hyperlink $ws.Range("A4")

# =============================================================================
# Sheet "de-de"
# =============================================================================
$ws = $wb.Worksheets.Item("de-de")

$ws.Rows.Item(3).Insert()
$ws.Hyperlinks.Delete()

# Row 2
$ws.Range("B2").Value2 = $statusReady
$ws.Range("D2").Value2 = $tsDe
$ws.Range("D2").NumberFormat = $dateFmt
$ws.Range("G2").Value2 = $epoch
$ws.Range("H2").Value2 = $reasonInc

# Row 3 (new)
$ws.Range("B3").Value2 = $statusReady
$ws.Range("D3").Value2 = $tsDe
$ws.Range("D3").NumberFormat = $dateFmt
$ws.Range("G3").Value2 = $epoch
$ws.Range("H3").Value2 = $reasonInc

# Row 4 (shifted ".localization-config" row)
$ws.Range("B4").Value2 = $statusNotLoc
$ws.Range("D4").Value2 = $epoch
$ws.Range("D4").NumberFormat = $dateFmt
$ws.Range("G4").Value2 = $epoch
$ws.Range("H4").Value2 = $reasonIgn

$ws.Hyperlinks.Add($ws.Range("A2"), ($mdUrlBase + $md1), "", "", $md1) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), ($xlfUrlBase + $xlfDe), "", "", $xlfDe) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), ($mdUrlBase + $md2), "", "", $md2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), ($xlfUrlBase + $xlfDe), "", "", $xlfDe) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), $cfgUrl, "", "", $cfgName) | Out-Null

Style-Hyperlink $ws.Range("A2")
Style-Hyperlink $ws.Range("C2")
Style-Hyperlink $ws.Range("A3")
Style-Hyperlink $ws.Range("C3")
Style-Hyperlink $ws.Range("A4")

Write-Output "Report regenerated for handoff."
